$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended at the bottom of the log (row 72).
# A72 looks like a date ("2025/10/07") but must stay literal text,
# matching the existing rows above it, so force text format before
# typing it in, then drop back to the default (unstyled) cell style.
$ws.Range("A72").NumberFormat = "@"
$ws.Range("A72").Value = "2025/10/07"
$ws.Range("A72").Style = "Normal"

$ws.Range("B72").Value = "火"
$ws.Range("C72").Value = 6
$ws.Range("D72").Value = 71
